$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 604.58826
$ws.Range("I28").Value = 526.7857
$ws.Range("J28").Value = 967.6667
$ws.Range("K28").Value = 526.7857
$ws.Range("L28").Value = 967.6667
$ws.Range("M28").Value = -41.78570000000002
$ws.Range("N28").Value = -1937.6667

$ws.Range("H46").Value = 2142.8572
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 2600
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 7800
$ws.Range("M46").Value = -2881
$ws.Range("N46").Value = -8038

$ws.Range("H60").Value = 2142.8572
$ws.Range("I60").Value = 1000
$ws.Range("J60").Value = 2600
$ws.Range("K60").Value = 3000
$ws.Range("L60").Value = 7800
$ws.Range("M60").Value = -2516
$ws.Range("N60").Value = -8768

$ws.Range("H113").Value = 3246.875
$ws.Range("I113").Value = 2670
$ws.Range("J113").Value = 3509.0908
$ws.Range("K113").Value = 2670
$ws.Range("L113").Value = 3509.0908
$ws.Range("M113").Value = 584
$ws.Range("N113").Value = -10017.0908

$ws.Range("H129").Value = 914.08826
$ws.Range("J129").Value = 1059.3334
$ws.Range("L129").Value = 3178.0002
$ws.Range("N129").Value = -13178.0002

$ws.Range("H138").Value = 2915.3125
$ws.Range("J138").Value = 3438.4856
$ws.Range("L138").Value = 10315.4568
$ws.Range("N138").Value = -20595.4568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 555.5263
$ws.Range("I2").Value = 499.27274
$ws.Range("K2").Value = 499.27274
$ws.Range("M2").Value = -386.27274

$ws.Range("H61").Value = 1263.159
$ws.Range("I61").Value = 1026.4595
$ws.Range("J61").Value = 2514.2856
$ws.Range("K61").Value = 1026.4595
$ws.Range("L61").Value = 2514.2856
$ws.Range("M61").Value = -814.4594999999999
$ws.Range("N61").Value = -2938.2856

$ws.Range("H97").Value = 2463.125
$ws.Range("I97").Value = 2943.923
$ws.Range("J97").Value = 379.66666
$ws.Range("K97").Value = 2943.923
$ws.Range("L97").Value = 379.66666
$ws.Range("M97").Value = -2447.923
$ws.Range("N97").Value = -1371.66666

$ws.Range("H116").Value = 555.5263
$ws.Range("I116").Value = 499.27274
$ws.Range("K116").Value = 499.27274
$ws.Range("M116").Value = 1794.72726

$ws.Range("H132").Value = 1972.5454
$ws.Range("I132").Value = 1486.238
$ws.Range("J132").Value = 2416.5652
$ws.Range("K132").Value = 4458.714
$ws.Range("L132").Value = 7249.6956
$ws.Range("M132").Value = -1928.714
$ws.Range("N132").Value = -12309.6956

$ws.Range("H136").Value = 1263.159
$ws.Range("I136").Value = 1026.4595
$ws.Range("J136").Value = 2514.2856
$ws.Range("K136").Value = 3079.3785
$ws.Range("L136").Value = 7542.8568
$ws.Range("M136").Value = -529.3784999999998
$ws.Range("N136").Value = -12642.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 555.5263
$ws.Range("I3").Value = 499.27274
$ws.Range("K3").Value = 499.27274
$ws.Range("M3").Value = -385.27274

$ws.Range("H134").Value = 2380.6365
$ws.Range("I134").Value = 1820.9546
$ws.Range("K134").Value = 5462.8638
$ws.Range("M134").Value = -2927.8638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5191.3184
$ws.Range("I134").Value = 5263.125
$ws.Range("K134").Value = 15789.375
$ws.Range("M134").Value = -13254.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 778.8
$ws.Range("I5").Value = 498.47058
$ws.Range("K5").Value = 1495.41174
$ws.Range("M5").Value = -1383.41174

$ws.Range("H80").Value = 951.8823
$ws.Range("I80").Value = 500.33334
$ws.Range("J80").Value = 1048.6428
$ws.Range("K80").Value = 1501.00002
$ws.Range("L80").Value = 3145.9284
$ws.Range("M80").Value = -565.0000199999999
$ws.Range("N80").Value = -5017.928400000001

$ws.Range("H83").Value = 951.8823
$ws.Range("I83").Value = 500.33334
$ws.Range("J83").Value = 1048.6428
$ws.Range("K83").Value = 4503.00006
$ws.Range("L83").Value = 9437.7852
$ws.Range("M83").Value = 176.9999399999997
$ws.Range("N83").Value = -18797.7852

$ws.Range("H86").Value = 347.5
$ws.Range("I86").Value = 347.5
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1042.5
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 143.5
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 347.5
$ws.Range("I89").Value = 347.5
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 3127.5
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 2800.5
$ws.Range("N89").ClearContents()

$ws.Range("H92").Value = 424.2857
$ws.Range("I92").Value = 405
$ws.Range("J92").Value = 472.5
$ws.Range("K92").Value = 1215
$ws.Range("L92").Value = 1417.5
$ws.Range("M92").Value = 33
$ws.Range("N92").Value = -3913.5

$ws.Range("H131").Value = 1961727
$ws.Range("J131").Value = 1013
$ws.Range("L131").Value = 3039
$ws.Range("N131").Value = -13119

$ws.Range("H135").Value = 778.8
$ws.Range("I135").Value = 498.47058
$ws.Range("K135").Value = 4486.23522
$ws.Range("M135").Value = -1951.23522

$ws.Range("H140").Value = 1277.5
$ws.Range("I140").Value = 866.72
$ws.Range("J140").Value = 2067.4614
$ws.Range("K140").Value = 2600.16
$ws.Range("L140").Value = 6202.3842
$ws.Range("M140").Value = 2579.84
$ws.Range("N140").Value = -16562.3842

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 13450
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 13450
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 13450
$ws.Range("N46").Value = -13762
$ws.Range("M46").ClearContents()

$ws.Range("H122").Value = 1378.2307
$ws.Range("I122").Value = 1523
$ws.Range("J122").Value = 1052.5
$ws.Range("K122").Value = 4569
$ws.Range("L122").Value = 3157.5
$ws.Range("M122").Value = -2119
$ws.Range("N122").Value = -8057.5

$ws.Range("H132").Value = 1860.8928
$ws.Range("J132").Value = 3542.4285
$ws.Range("L132").Value = 10627.2855
$ws.Range("N132").Value = -15687.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3285.3333
$ws.Range("I61").Value = 3285.3333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3285.3333
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3083.3333
$ws.Range("N61").ClearContents()

$ws.Range("H113").Value = 3285.3333
$ws.Range("I113").Value = 3285.3333
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3285.3333
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1115.3333
$ws.Range("N113").ClearContents()

$ws.Range("H136").Value = 4833.3613
$ws.Range("I136").Value = 7325.8
$ws.Range("J136").Value = 1717.8125
$ws.Range("K136").Value = 21977.4
$ws.Range("L136").Value = 5153.4375
$ws.Range("M136").Value = -19427.4
$ws.Range("N136").Value = -10253.4375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 47621210
$ws.Range("I81").Value = 100002216
$ws.Range("J81").Value = 2109.182
$ws.Range("K81").Value = 200004432
$ws.Range("L81").Value = 4218.364
$ws.Range("M81").Value = -200003371
$ws.Range("N81").Value = -6340.364

$ws.Range("H84").Value = 47621210
$ws.Range("I84").Value = 100002216
$ws.Range("J84").Value = 2109.182
$ws.Range("K84").Value = 1000022160
$ws.Range("L84").Value = 21091.82
$ws.Range("M84").Value = -1000016856
$ws.Range("N84").Value = -31699.82

$ws.Range("H113").Value = 29411974
$ws.Range("I113").Value = 222.8125
$ws.Range("J113").Value = 500000000
$ws.Range("K113").Value = 668.4375
$ws.Range("L113").Value = 1500000000
$ws.Range("M113").Value = 1501.5625
$ws.Range("N113").Value = -1500004340

$ws.Range("H136").Value = 837.72095
$ws.Range("I136").Value = 540.0833
$ws.Range("J136").Value = 2368.4285
$ws.Range("K136").Value = 1620.2499
$ws.Range("L136").Value = 7105.2855
$ws.Range("M136").Value = 929.7501
$ws.Range("N136").Value = -12205.2855
